$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert "~30" labels to plain numbers ---
$ws.Range("C2").Value = 30
$ws.Range("C3").Value = 30
$ws.Range("C4").Value = 30
$ws.Range("C5").Value = 30
$ws.Range("C6").Value = 30

# --- Row 13 ListView: "3 hours" -> 180 (minutes) ---
$ws.Range("D13").Value = 180

# --- Row 14: MenuBar -> GarageMenuBar ---
$ws.Range("A14").Value = "GarageMenuBar"

# --- Duration strings -> plain minute numbers ---
$ws.Range("D15").Value = 60
$ws.Range("D16").Value = 30
$ws.Range("D17").Value = 20
$ws.Range("D18").Value = 30
$ws.Range("D19").Value = 20

# --- Row 22 DataSaver: drop B22 (20), convert D22 "20 minutes" -> 20 ---
$ws.Range("B22").ClearContents()
$ws.Range("D22").Value = 20

# --- Row 24 Javadoc: "~100" -> 100 ---
$ws.Range("C24").Value = 100

# --- Row 25: "120 (General DB)" -> 120 ---
$ws.Range("C25").Value = 120

# --- Row 27 LightWork: give it the same green highlight as the other task names ---
$ws.Range("A27").Interior.Color = 5287936

# --- Clear old totals row 28 (will be replaced by new layout ending at row 34) ---
$ws.Range("B28").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("D28").ClearContents()

# --- New tasks appended (rows 28-33), styled like the other green task-name cells ---
$ws.Range("A28").Value = "LoginButton"
$ws.Range("B28").Value = 30
$ws.Range("A29").Value = "ViewGarageButton"
$ws.Range("B29").Value = 30
$ws.Range("A30").Value = "DeleteUserDialog"
$ws.Range("B30").Value = 40
$ws.Range("A31").Value = "MainMenuPane"
$ws.Range("B31").Value = 10
$ws.Range("A32").Value = "MainStage"
$ws.Range("B32").Value = 10
$ws.Range("A33").Value = "UserDataViewer"
$ws.Range("B33").Value = 180

$ws.Range("A28:A33").Interior.Color = 5287936
$ws.Range("A28:A33").HorizontalAlignment = -4108

# --- New totals row 34 ---
$ws.Range("B34").Formula = "=SUM(B2:B33)"
$ws.Range("C34").Formula = "=SUM(C2:C33)"
$ws.Range("D34").Formula = "=SUM(D2:D33)"
$ws.Range("E34").Value = "<---- (Minutes)"

# --- Header A1 gets an underline ---
$ws.Range("A1").Font.Underline = $true

# --- Selection / scroll state to match the edited document ---
$ws.Range("E33").Select()
